$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 273
$ws.Range("I2").Value = 202.25
$ws.Range("J2").Value = 329.6
$ws.Range("K2").Value = 202.25
$ws.Range("L2").Value = 329.6
$ws.Range("M2").Value = -89.25
$ws.Range("N2").Value = -555.6
$ws.Range("H6").Value = 1110
$ws.Range("I6").Value = 1243.2
$ws.Range("K6").Value = 3729.6
$ws.Range("M6").Value = -3617.6
$ws.Range("H55").Value = 400
$ws.Range("I55").Value = 400
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 400
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = -186
$ws.Range("N55").ClearContents()
$ws.Range("H141").Value = 3559.4285
$ws.Range("I141").Value = 3314.875
$ws.Range("J141").Value = 3885.5
$ws.Range("K141").Value = 9944.625
$ws.Range("L141").Value = 11656.5
$ws.Range("M141").Value = -4764.625
$ws.Range("N141").Value = -22016.5

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 863522.0600000001
$ws.Range("I2").Value = 1369824.4
$ws.Range("J2").Value = 2808
$ws.Range("K2").Value = 1369824.4
$ws.Range("L2").Value = 2808
$ws.Range("M2").Value = -1369711.4
$ws.Range("N2").Value = -3034
$ws.Range("H32").Value = 4489.6924
$ws.Range("I32").Value = 3551.7144
$ws.Range("J32").Value = 8429.200000000001
$ws.Range("K32").Value = 3551.7144
$ws.Range("L32").Value = 8429.200000000001
$ws.Range("M32").Value = -3264.7144
$ws.Range("N32").Value = -9003.200000000001
$ws.Range("H45").Value = 7501607.5
$ws.Range("J45").Value = 2250
$ws.Range("L45").Value = 2250
$ws.Range("N45").Value = -3004
$ws.Range("H116").Value = 863522.0600000001
$ws.Range("I116").Value = 1369824.4
$ws.Range("J116").Value = 2808
$ws.Range("K116").Value = 1369824.4
$ws.Range("L116").Value = 2808
$ws.Range("M116").Value = -1367530.4
$ws.Range("N116").Value = -7396
$ws.Range("H122").Value = 1475.9166
$ws.Range("I122").Value = 1559.1578
$ws.Range("K122").Value = 4677.4734
$ws.Range("M122").Value = -2227.4734
$ws.Range("H132").Value = 1745.5745
$ws.Range("I132").Value = 1313.3103
$ws.Range("K132").Value = 3939.9309
$ws.Range("M132").Value = -1409.9309

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 863522.0600000001
$ws.Range("I3").Value = 1369824.4
$ws.Range("J3").Value = 2808
$ws.Range("K3").Value = 1369824.4
$ws.Range("L3").Value = 2808
$ws.Range("M3").Value = -1369710.4
$ws.Range("N3").Value = -3036
$ws.Range("H107").Value = 4321.5557
$ws.Range("I107").Value = 2225
$ws.Range("K107").Value = 2225
$ws.Range("M107").Value = -305
$ws.Range("H134").Value = 5138.886
$ws.Range("J134").Value = 3181.111
$ws.Range("L134").Value = 9543.332999999999
$ws.Range("N134").Value = -14613.333

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 199.5
$ws.Range("I7").Value = 199.66667
$ws.Range("J7").Value = 199
$ws.Range("K7").Value = 199.66667
$ws.Range("L7").Value = 199
$ws.Range("M7").Value = -86.66667000000001
$ws.Range("N7").Value = -425
$ws.Range("H22").Value = 1555.3334
$ws.Range("I22").Value = 1166.5
$ws.Range("K22").Value = 1166.5
$ws.Range("M22").Value = -816.5
$ws.Range("H31").Value = 1730.4231
$ws.Range("J31").Value = 2534.6
$ws.Range("L31").Value = 2534.6
$ws.Range("N31").Value = -3124.6
$ws.Range("H34").Value = 1730.4231
$ws.Range("J34").Value = 2534.6
$ws.Range("L34").Value = 2534.6
$ws.Range("N34").Value = -2938.6
$ws.Range("H99").Value = 5000949.5
$ws.Range("J99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("N99").ClearContents()
$ws.Range("H126").Value = 5000949.5
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 303.74075
$ws.Range("I5").Value = 247.07692
$ws.Range("K5").Value = 741.23076
$ws.Range("M5").Value = -629.23076
$ws.Range("H51").Value = 2099.8
$ws.Range("J51").Value = 2374.75
$ws.Range("L51").Value = 7124.25
$ws.Range("N51").Value = -8044.25
$ws.Range("H68").Value = 2429
$ws.Range("I68").Value = 1266
$ws.Range("J68").Value = 2910.2415
$ws.Range("K68").Value = 3798
$ws.Range("L68").Value = 8730.7245
$ws.Range("M68").Value = -2987
$ws.Range("N68").Value = -10352.7245
$ws.Range("H71").Value = 2429
$ws.Range("I71").Value = 1266
$ws.Range("J71").Value = 2910.2415
$ws.Range("K71").Value = 11394
$ws.Range("L71").Value = 26192.1735
$ws.Range("M71").Value = -7338
$ws.Range("N71").Value = -34304.1735
$ws.Range("H131").Value = 12213338
$ws.Range("J131").Value = 19123.871
$ws.Range("L131").Value = 57371.613
$ws.Range("N131").Value = -67451.613
$ws.Range("H133").Value = 86333400
$ws.Range("I133").Value = 86333400
$ws.Range("K133").Value = 259000200
$ws.Range("M133").Value = -258995140
$ws.Range("H135").Value = 303.74075
$ws.Range("I135").Value = 247.07692
$ws.Range("K135").Value = 2223.69228
$ws.Range("M135").Value = 311.3077199999998
$ws.Range("H137").Value = 3092.24
$ws.Range("I137").Value = 786
$ws.Range("J137").Value = 4904.2856
$ws.Range("K137").Value = 2358
$ws.Range("L137").Value = 14712.8568
$ws.Range("M137").Value = 2742
$ws.Range("N137").Value = -24912.8568

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 21271.428
$ws.Range("J46").Value = 21271.428
$ws.Range("L46").Value = 21271.428
$ws.Range("N46").Value = -21583.428
$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()
$ws.Range("H70").Value = 20280
$ws.Range("I70").Value = 30666.666
$ws.Range("J70").Value = 4700
$ws.Range("K70").Value = 30666.666
$ws.Range("L70").Value = 4700
$ws.Range("M70").Value = -30396.666
$ws.Range("N70").Value = -5240
$ws.Range("H73").Value = 20280
$ws.Range("I73").Value = 30666.666
$ws.Range("J73").Value = 4700
$ws.Range("K73").Value = 30666.666
$ws.Range("L73").Value = 4700
$ws.Range("M73").Value = -29730.666
$ws.Range("N73").Value = -6572
$ws.Range("H113").Value = 1700
$ws.Range("J113").Value = 1700
$ws.Range("L113").Value = 1700
$ws.Range("N113").Value = -6040
$ws.Range("H122").Value = 2146.9583
$ws.Range("I122").Value = 1567.2941
$ws.Range("K122").Value = 4701.8823
$ws.Range("M122").Value = -2251.8823

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 6819.2856
$ws.Range("I16").Value = 7861
$ws.Range("K16").Value = 7861
$ws.Range("M16").Value = -7691
$ws.Range("H22").Value = 2017.5454
$ws.Range("I22").Value = 1732.6666
$ws.Range("J22").Value = 2124.375
$ws.Range("K22").Value = 1732.6666
$ws.Range("L22").Value = 2124.375
$ws.Range("M22").Value = -1437.6666
$ws.Range("N22").Value = -2714.375
$ws.Range("H27").Value = 2017.5454
$ws.Range("I27").Value = 1732.6666
$ws.Range("J27").Value = 2124.375
$ws.Range("K27").Value = 1732.6666
$ws.Range("L27").Value = 2124.375
$ws.Range("M27").Value = -1625.6666
$ws.Range("N27").Value = -2338.375
$ws.Range("H46").Value = 1303.7333
$ws.Range("I46").Value = 474
$ws.Range("J46").Value = 1431.3846
$ws.Range("K46").Value = 474
$ws.Range("L46").Value = 1431.3846
$ws.Range("M46").Value = -286
$ws.Range("N46").Value = -1807.3846
$ws.Range("H82").Value = 3982.5715
$ws.Range("I82").Value = 1600
$ws.Range("J82").Value = 4935.6
$ws.Range("K82").Value = 1600
$ws.Range("L82").Value = 4935.6
$ws.Range("M82").Value = -1239
$ws.Range("N82").Value = -5657.6
$ws.Range("H85").Value = 3982.5715
$ws.Range("I85").Value = 1600
$ws.Range("J85").Value = 4935.6
$ws.Range("K85").Value = 1600
$ws.Range("L85").Value = 4935.6
$ws.Range("M85").Value = -352
$ws.Range("N85").Value = -7431.6
$ws.Range("H100").Value = 4175.6
$ws.Range("I100").Value = 3949.5
$ws.Range("K100").Value = 3949.5
$ws.Range("M100").Value = -3408.5

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 3249.75
$ws.Range("J96").Value = 3249.75
$ws.Range("L96").Value = 3249.75
$ws.Range("N96").Value = -5995.75
